$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 265.0625
$ws.Range("I39").Value = 17.357143
$ws.Range("K39").Value = 52.071429
$ws.Range("M39").Value = 243.928571

$ws.Range("H43").Value = 11189.8
$ws.Range("J43").Value = 11189.8
$ws.Range("L43").Value = 11189.8
$ws.Range("N43").Value = -11327.8

$ws.Range("H76").Value = 4900
$ws.Range("J76").Value = 10000
$ws.Range("L76").Value = 10000
$ws.Range("N76").Value = -10630

$ws.Range("H79").Value = 4900
$ws.Range("J79").Value = 10000
$ws.Range("L79").Value = 10000
$ws.Range("N79").Value = -12184

$ws.Range("H137").Value = 2320.6924
$ws.Range("I137").Value = 2021.95
$ws.Range("K137").Value = 6065.85
$ws.Range("M137").Value = -3515.85

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6442.077
$ws.Range("I61").Value = 6442.077
$ws.Range("K61").Value = 6442.077
$ws.Range("M61").Value = -6230.077

$ws.Range("H97").Value = 294.86667
$ws.Range("I97").Value = 286.3846
$ws.Range("K97").Value = 286.3846
$ws.Range("M97").Value = 209.6154

$ws.Range("H102").Value = 2156.8147
$ws.Range("I102").Value = 1855.1538
$ws.Range("K102").Value = 1855.1538
$ws.Range("M102").Value = -233.1538

$ws.Range("H136").Value = 6442.077
$ws.Range("I136").Value = 6442.077
$ws.Range("K136").Value = 19326.231
$ws.Range("M136").Value = -16776.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1560.1052
$ws.Range("I20").Value = 1503.0714
$ws.Range("K20").Value = 1503.0714
$ws.Range("M20").Value = -1256.0714

$ws.Range("H86").Value = 19233878
$ws.Range("I86").Value = 2208.913
$ws.Range("K86").Value = 2208.913
$ws.Range("M86").Value = -1085.913

$ws.Range("H89").Value = 19233878
$ws.Range("I89").Value = 2208.913
$ws.Range("K89").Value = 11044.565
$ws.Range("M89").Value = -5428.565000000001

$ws.Range("H94").Value = 2668.45
$ws.Range("I94").Value = 2282.5789
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 2282.5789
$ws.Range("L94").Value = 10000
$ws.Range("M94").Value = -1831.5789
$ws.Range("N94").Value = -10902

$ws.Range("H132").Value = 63748
$ws.Range("J132").Value = 63748
$ws.Range("L132").Value = 63748
$ws.Range("N132").Value = -73868

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 2071.6
$ws.Range("I13").Value = 4
$ws.Range("J13").Value = 3450
$ws.Range("K13").Value = 4
$ws.Range("L13").Value = 3450
$ws.Range("M13").Value = 135
$ws.Range("N13").Value = -3728

$ws.Range("H58").Value = 2026.8422
$ws.Range("I58").Value = 2026.8422
$ws.Range("K58").Value = 2026.8422
$ws.Range("M58").Value = -1823.8422

$ws.Range("H132").Value = 3549.8462
$ws.Range("I132").Value = 3549.8462
$ws.Range("K132").Value = 10649.5386
$ws.Range("M132").Value = -8119.5386

$ws.Range("H136").Value = 2026.8422
$ws.Range("I136").Value = 2026.8422
$ws.Range("K136").Value = 6080.5266
$ws.Range("M136").Value = -3530.5266

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11162604
$ws.Range("I4").Value = 10241776
$ws.Range("J4").Value = 30500000
$ws.Range("K4").Value = 30725328
$ws.Range("L4").Value = 91500000
$ws.Range("M4").Value = -30725216
$ws.Range("N4").Value = -91500224

$ws.Range("H16").Value = 231.14285
$ws.Range("I16").Value = 203
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 609
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -436
$ws.Range("N16").Value = -1546

$ws.Range("H28").Value = 1309.8
$ws.Range("I28").Value = 1309.8
$ws.Range("K28").Value = 3929.4
$ws.Range("M28").Value = -3697.4

$ws.Range("H92").Value = 4874.25
$ws.Range("I92").Value = 6249
$ws.Range("K92").Value = 18747
$ws.Range("M92").Value = -17499

$ws.Range("H107").Value = 2467.3333
$ws.Range("I107").Value = 4158.7144
$ws.Range("K107").Value = 12476.1432
$ws.Range("M107").Value = -10556.1432

$ws.Range("H131").Value = 1589480
$ws.Range("J131").Value = 2022788.1
$ws.Range("L131").Value = 6068364.300000001
$ws.Range("N131").Value = -6078444.300000001

$ws.Range("H136").Value = 1907.5
$ws.Range("I136").Value = 1907.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5722.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -622.5
$ws.Range("N136").ClearContents()

$ws.Range("H140").Value = 1988.1428
$ws.Range("I140").Value = 1320
$ws.Range("K140").Value = 3960
$ws.Range("M140").Value = 1220

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 522.8823
$ws.Range("I2").Value = 54.333332
$ws.Range("K2").Value = 54.333332
$ws.Range("M2").Value = 58.666668

$ws.Range("H70").Value = 27784546
$ws.Range("J70").Value = 7689.778
$ws.Range("L70").Value = 7689.778
$ws.Range("N70").Value = -8229.778

$ws.Range("H73").Value = 27784546
$ws.Range("J73").Value = 7689.778
$ws.Range("L73").Value = 7689.778
$ws.Range("N73").Value = -9561.778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 3063
$ws.Range("J12").Value = 3063
$ws.Range("L12").Value = 3063
$ws.Range("N12").Value = -3403

$ws.Range("H16").Value = 801.3333
$ws.Range("I16").Value = 561.8
$ws.Range("K16").Value = 561.8
$ws.Range("M16").Value = -391.8

$ws.Range("H46").Value = 983
$ws.Range("I46").Value = 979.8
$ws.Range("J46").Value = 999
$ws.Range("K46").Value = 979.8
$ws.Range("L46").Value = 999
$ws.Range("M46").Value = -791.8
$ws.Range("N46").Value = -1375

$ws.Range("H55").Value = 1335.9231
$ws.Range("I55").Value = 675.8823
$ws.Range("K55").Value = 675.8823
$ws.Range("M55").Value = -502.8823

$ws.Range("H68").Value = 8840.583000000001
$ws.Range("I68").Value = 2759.4
$ws.Range("K68").Value = 2759.4
$ws.Range("M68").Value = -2010.4

$ws.Range("H71").Value = 8840.583000000001
$ws.Range("I71").Value = 2759.4
$ws.Range("K71").Value = 13797
$ws.Range("M71").Value = -10053

$ws.Range("H100").Value = 10400
$ws.Range("I100").Value = 4000
$ws.Range("J100").Value = 11111.111
$ws.Range("K100").Value = 4000
$ws.Range("L100").Value = 11111.111
$ws.Range("M100").Value = -3459
$ws.Range("N100").Value = -12193.111

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1992.6316
$ws.Range("I81").Value = 1521.1765
$ws.Range("J81").Value = 6000
$ws.Range("K81").Value = 3042.353
$ws.Range("L81").Value = 12000
$ws.Range("M81").Value = -1981.353
$ws.Range("N81").Value = -14122

$ws.Range("H84").Value = 1992.6316
$ws.Range("I84").Value = 1521.1765
$ws.Range("J84").Value = 6000
$ws.Range("K84").Value = 15211.765
$ws.Range("L84").Value = 60000
$ws.Range("M84").Value = -9907.764999999999
$ws.Range("N84").Value = -70608

$ws.Range("H107").Value = 1548.8
$ws.Range("I107").Value = 1063
$ws.Range("J107").Value = 1757
$ws.Range("K107").Value = 3189
$ws.Range("L107").Value = 5271
$ws.Range("M107").Value = -1269
$ws.Range("N107").Value = -9111

$ws.Range("H136").Value = 1559.2903
$ws.Range("I136").Value = 1383.3334
$ws.Range("J136").Value = 2747
$ws.Range("K136").Value = 4150.0002
$ws.Range("L136").Value = 8241
$ws.Range("M136").Value = -1600.0002
$ws.Range("N136").Value = -13341
